$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "46.923.94"
$ws.Range("E2").Value = "  +4.32%  "
$ws.Range("D3").Value = "2.265.61"
$ws.Range("E3").Value = "  -0.31%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'300.86"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.41%  "
$ws.Range("D6").Value = "'100.17"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.11%  "
$ws.Range("E7").Value = "  -0.64%  "
$ws.Range("D8").Value = "'0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "'0.512"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Value = "'35.64"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.99%  "
$ws.Range("E11").Value = "  -1.22%  "
$ws.Range("D12").Value = "'7.17"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.67%  "
$ws.Range("E13").Value = "  -0.89%  "
$ws.Range("D14").Value = "2.609.87"
$ws.Range("E14").Value = "  -0.13%  "
$ws.Range("D15").Value = "2.263.42"
$ws.Range("E15").Value = "  -0.11%  "
$ws.Range("D16").Value = "'13.58"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "46.846.73"
$ws.Range("E17").Value = "  +4.47%  "
$ws.Range("D18").Value = "'0.795"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.90%  "
$ws.Range("D19").Value = "'12.76"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.36%  "
$ws.Range("D20").Value = "0.0₃0926"
$ws.Range("E20").Value = "  +0.51%  "
$ws.Range("D21").Value = "'5.86"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.44%  "
$ws.Range("D22").Value = "'65.10"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.93%  "
$ws.Range("D23").Value = "'248.94"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.33%  "
$ws.Range("D24").Value = "'2.81"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.57%  "
$ws.Range("D25").Value = "'0.999"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.15%  "
$ws.Range("E26").Value = "  -1.07%  "
$ws.Range("D27").Value = "'42.29"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.11%  "
$ws.Range("E28").Value = "  -2.59%  "
$ws.Range("D29").Value = "'9.71"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.89%  "
$ws.Range("D30").Value = "'19.86"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.13%  "
$ws.Range("D31").Value = "'2.79"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +8.91%  "
$ws.Range("D32").Value = "'145.35"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.67%  "
$ws.Range("D33").Value = "'5.41"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.81%  "
$ws.Range("D34").Value = "'0.0775"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.22%  "
$ws.Range("E35").Value = "  +8.12%  "
$ws.Range("D36").Value = "'0.115"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +10.31%  "
$ws.Range("E37").Value = "  -1.62%  "
$ws.Range("D38").Value = "'16.25"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +18.45%  "
$ws.Range("E39").Value = "  -3.82%  "
$ws.Range("D40").Value = "'3.80"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.53%  "
$ws.Range("D41").Value = "'0.0298"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.77%  "
$ws.Range("D42").Value = "'3.21"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.87%  "
$ws.Range("D43").Value = "'0.999"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.09%  "
$ws.Range("D44").Value = "'1.97"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.34%  "
$ws.Range("B45").Value = "BitcoinSV"
$ws.Range("C45").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D45").Value = "'91.41"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +18.78%  "
$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "1.783.51"
$ws.Range("E46").Value = "  +2.35%  "
$ws.Range("E47").Value = "  -3.45%  "
$ws.Range("D48").Value = "'71.41"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.58%  "
$ws.Range("D49").Value = "'4.82"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.77%  "
$ws.Range("D50").Value = "'93.89"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.82%  "
$ws.Range("D51").Value = "'7.82"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.67%  "
